$wb = $excel.ActiveWorkbook

# The original "Hoja1" sheet becomes "TypesWithErrors" (kept in place, same
# data) and a brand new "TypesOK" sheet (clean data, no errors) is inserted
# before it so the final tab order is TypesOK, TypesWithErrors.
$errors = $wb.Worksheets.Item(1)
$errors.Name = "TypesWithErrors"

# Fix the header typo on the existing sheet ("IntColumn" -> "Int Column").
$errors.Range("A1").Value = "Int Column"

$ok = $wb.Worksheets.Add($errors)
$ok.Name = "TypesOK"

# Headers
$ok.Range("A1").Value = "Int Column"
$ok.Range("B1").Value = "String Column"
$ok.Range("C1").Value = "Date Column"
$ok.Range("D1").Value = "Bool column"

# Row 2
$ok.Range("A2").Value = 1
$ok.Range("B2").Value = "Item 1"
$ok.Range("C2").Value = [DateTime]"2000-01-01"
$ok.Range("D2").Value = 1

# Row 3
$ok.Range("A3").Value = 2
$ok.Range("B3").Value = "Item 2"
$ok.Range("C3").Value = [DateTime]"2000-01-02"
$ok.Range("D3").Value = "Y"

# Row 4
$ok.Range("A4").Value = 3
$ok.Range("B4").Value = "Item 3"
$ok.Range("C4").Value = [DateTime]"2000-01-03"
$ok.Range("D4").Value = 0

# Row 5
$ok.Range("A5").Value = 5
$ok.Range("B5").Value = "Item 4"
$ok.Range("C5").Value = [DateTime]"2000-01-04"
$ok.Range("D5").Value = "N"

# Row 6
$ok.Range("A6").Value = 6
$ok.Range("B6").Value = "Item 5"
$ok.Range("C6").Value = [DateTime]"2000-01-05"
$ok.Range("D6").Value = "S"

[void]$ok.Range("D4").Select()
